$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioB")
$ws.Range("D17").Value = 1.5
$ws.Range("D18").Value = 0.5
$ws.Range("D24").Value = 1.5
$ws.Range("D25").Value = 0.5
